$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D:E columns to text format so numeric-looking strings (e.g. "1.00",
# "34.456.73", percentages) are stored as literal text, matching the source
# workbook which stores these as inline strings, not numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "34.456.73"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "1.801.01"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "224.57"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").Value = "0.603"
$ws.Range("E6").Value = "  +1.68%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").Value = "41.23"
$ws.Range("E8").Value = "  +14.04%  "
$ws.Range("B9").Value = "WrappedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D9").Value = "3.173.00"
$ws.Range("E9").Value = "  +75.63%  "
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "0.290"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "0.0668"
$ws.Range("E11").Value = "  -2.42%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.0991"
$ws.Range("E12").Value = "  +2.83%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.059.66"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "10.83"
$ws.Range("E14").Value = "  -3.40%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.630"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "34.466.89"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("D17").Value = "4.38"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").Value = "67.39"
$ws.Range("E18").Value = "  -3.17%  "
$ws.Range("D19").Value = "240.13"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D20").Value = "0.0₃0765"
$ws.Range("E20").Value = "  -2.03%  "
$ws.Range("D21").Value = "11.07"
$ws.Range("E21").Value = "  -2.27%  "
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").Value = "4.08"
$ws.Range("E23").Value = "  -1.21%  "
$ws.Range("E24").Value = "  -3.98%  "
$ws.Range("D25").Value = "171.78"
$ws.Range("D26").Value = "7.66"
$ws.Range("E26").Value = "  -3.48%  "
$ws.Range("D27").Value = "17.36"
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("D28").Value = "0.121"
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("D31").Value = "3.78"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("D33").Value = "3.84"
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("D34").Value = "1.78"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.644"
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "1.314.17"
$ws.Range("E36").Value = "  -3.46%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("D38").Value = "14.95"
$ws.Range("E38").Value = "  +12.41%  "
$ws.Range("B39").Value = "Aave"
$ws.Range("C39").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D39").Value = "84.89"
$ws.Range("E39").Value = "  +4.04%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.0187"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "1.25"
$ws.Range("E41").Value = "  +7.46%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "2.33"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").Value = "2.43"
$ws.Range("E43").Value = "  +0.53%  "
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").Value = "0.937"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("E46").Value = "  +4.47%  "
$ws.Range("D47").Value = "1.960.38"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").Value = "5.76"
$ws.Range("E48").Value = "  -1.36%  "
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").Value = "100.57"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("D51").Value = "0.0610"
$ws.Range("E51").Value = "  +1.01%  "

# Restore the default (Normal) style so no stray number-format style sticks
# to the cells (matches the unstyled inlineStr cells in the target).
$ws.Range("D2:E51").Style = "Normal"
